# Rotate the "Recorded By" list in column G: move the first (leftmost)
# comma-separated entry to the end of the list, for every row where more
# than one entry is present. Single-entry cells are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value2 = $rotated
    }
}
